$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "SebastianJerez"
$ws.Range("B6").Value = "Starjerez"
$ws.Range("C6").Value = "27242679jsjs"
$ws.Range("D6").Value = "sebastianjs99@hotmail.com"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 5
